# Trade #14 closed at 2026-02-17 15:17:41 - unknown UNKNOWN +0.000%
#
# - Updates the Summary sheet totals (Current Capital, Total P&L $/%%,
#   Total/Winning Trades, Win Rate %%) to reflect the newly closed trade.
# - Updates the MarketMaking strategy row on the "Strategy Status" sheet
#   to match.
# - Appends the new trade (#14) as row 15 on both the "All Trades" sheet
#   and the per-strategy "MarketMaking" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.8    # Current Capital
$summary.Range("B4").Value = -0.2      # Total P&L $
$summary.Range("B5").Value = -0.29     # Total P&L %
$summary.Range("B6").Value = 14        # Total Trades
$summary.Range("B7").Value = 4         # Winning Trades
$summary.Range("B9").Value = 28.57     # Win Rate %

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.8       # Capital
$status.Range("D4").Value = 14         # Trades
$status.Range("E4").Value = -0.2       # P&L $
$status.Range("F4").Value = -0.2       # P&L %
$status.Range("G4").Value = 28.57      # Win Rate %

# ---------------------------------------------------------------
# 3) Append trade #14 to a trade-log sheet (shared helper)
# ---------------------------------------------------------------
function Add-Trade14Row($ws) {
    $row = 15
    $ws.Cells.Item($row, 1).Value = 14
    # The date string would otherwise be auto-detected as a date by the
    # COM layer and stored as a date serial - force plain text with a
    # leading apostrophe (like typing it in the Excel UI), then clear the
    # resulting quote-prefix style so the cell is left unstyled.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Value = "'2026-02-17"
    $dateCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "15:17:34"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.8
    $ws.Cells.Item($row, 7).Value = 0.85
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 6.25
    $ws.Cells.Item($row, 10).Value = 0.05
    $ws.Cells.Item($row, 11).Value = 99.8
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.15
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade14Row $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade14Row $marketMaking
